# Generate Report for Handoff
# Adds a new row (for file 86c14bdc-2d76-44ad-bb1b-cbe6d32d2268.md) to each
# of the three report tables: Overview, zh-cn, de-de.

$wb = $excel.ActiveWorkbook

$commitSha = "ad83d812831d97b9452e9c2d81804c936654ff90"
$newGuid   = "86c14bdc-2d76-44ad-bb1b-cbe6d32d2268"
$newMd     = $newGuid + ".md"
$newMdPath = "e2e\" + $newMd
$newUrl    = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/" + $commitSha + "/e2e/" + $newMd

$dateFmt = "yyyy-mm-dd HH:mm:ss"

# ---------------------------------------------------------------------
# Sheet "Overview" -- one summary row per handed-off source file
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$loOverview = $wsOverview.ListObjects.Item(1)
$newOvRow = $loOverview.ListRows.Add()

$ovRow = $newOvRow.Range.Row

$wsOverview.Cells.Item($ovRow, 1).Value = $newMd
$wsOverview.Cells.Item($ovRow, 2).Value = $newMdPath
$wsOverview.Cells.Item($ovRow, 3).Value = ".md"
$wsOverview.Cells.Item($ovRow, 4).Value = ""
$wsOverview.Cells.Item($ovRow, 5).Value = "Ready for handoff"
$wsOverview.Cells.Item($ovRow, 6).Value = "Ready for handoff"
$wsOverview.Cells.Item($ovRow, 7).Value = "2016-09-08 04:53:57"
$wsOverview.Cells.Item($ovRow, 7).NumberFormat = $dateFmt

$wsOverview.Hyperlinks.Add($wsOverview.Cells.Item($ovRow, 2), $newUrl, [System.Type]::Missing, [System.Type]::Missing, $newMdPath) | Out-Null

# ---------------------------------------------------------------------
# Sheet "zh-cn" -- per-language handoff detail row
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$loZh = $wsZh.ListObjects.Item(1)
$newZhRow = $loZh.ListRows.Add()

$zhRow = $newZhRow.Range.Row

$zhXlf = $newGuid + ".fed3b682a84d7756a9941ccca5215e9e3f67e2a4.zh-cn.xlf"

$wsZh.Cells.Item($zhRow, 1).Value = $newMd
$wsZh.Cells.Item($zhRow, 2).Value = ".md"
$wsZh.Cells.Item($zhRow, 3).Value = "Ready for handoff"
$wsZh.Cells.Item($zhRow, 4).Value = "e2e"
$wsZh.Cells.Item($zhRow, 5).Value = "ht"
$wsZh.Cells.Item($zhRow, 6).Value = "False"
$wsZh.Cells.Item($zhRow, 7).Value = $zhXlf
$wsZh.Cells.Item($zhRow, 8).Value = "2016-09-08 04:53:51"
$wsZh.Cells.Item($zhRow, 8).NumberFormat = $dateFmt
$wsZh.Cells.Item($zhRow, 9).Value = ""
$wsZh.Cells.Item($zhRow, 10).Value = ""
$wsZh.Cells.Item($zhRow, 11).Value = "0001-01-01 00:00:00"
$wsZh.Cells.Item($zhRow, 11).NumberFormat = $dateFmt
$wsZh.Cells.Item($zhRow, 12).Value = ""
$wsZh.Cells.Item($zhRow, 13).Value = "True"
$wsZh.Cells.Item($zhRow, 14).Value = ""
$wsZh.Cells.Item($zhRow, 15).Value = "False"
$wsZh.Cells.Item($zhRow, 16).Value = ""

$wsZh.Hyperlinks.Add($wsZh.Cells.Item($zhRow, 1), $newUrl, [System.Type]::Missing, [System.Type]::Missing, $newMd) | Out-Null

# ---------------------------------------------------------------------
# Sheet "de-de" -- per-language handoff detail row
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$loDe = $wsDe.ListObjects.Item(1)
$newDeRow = $loDe.ListRows.Add()

$deRow = $newDeRow.Range.Row

$deXlf = $newGuid + ".fed3b682a84d7756a9941ccca5215e9e3f67e2a4.de-de.xlf"

$wsDe.Cells.Item($deRow, 1).Value = $newMd
$wsDe.Cells.Item($deRow, 2).Value = ".md"
$wsDe.Cells.Item($deRow, 3).Value = "Ready for handoff"
$wsDe.Cells.Item($deRow, 4).Value = "e2e"
$wsDe.Cells.Item($deRow, 5).Value = "ht"
$wsDe.Cells.Item($deRow, 6).Value = "False"
$wsDe.Cells.Item($deRow, 7).Value = $deXlf
$wsDe.Cells.Item($deRow, 8).Value = "2016-09-08 04:53:57"
$wsDe.Cells.Item($deRow, 8).NumberFormat = $dateFmt
$wsDe.Cells.Item($deRow, 9).Value = ""
$wsDe.Cells.Item($deRow, 10).Value = ""
$wsDe.Cells.Item($deRow, 11).Value = "0001-01-01 00:00:00"
$wsDe.Cells.Item($deRow, 11).NumberFormat = $dateFmt
$wsDe.Cells.Item($deRow, 12).Value = ""
$wsDe.Cells.Item($deRow, 13).Value = "True"
$wsDe.Cells.Item($deRow, 14).Value = ""
$wsDe.Cells.Item($deRow, 15).Value = "False"
$wsDe.Cells.Item($deRow, 16).Value = ""

$wsDe.Hyperlinks.Add($wsDe.Cells.Item($deRow, 1), $newUrl, [System.Type]::Missing, [System.Type]::Missing, $newMd) | Out-Null

Write-Host "Added handoff rows for $newMd to Overview, zh-cn, de-de"
